# Generate Report for Handback
#
# The handback transform failed for the "26c9a277-901b-4ea2-b1f8-1505ff629220"
# source file in both the zh-cn and de-de locale sheets. This updates the
# Status everywhere it is shown (the per-locale sheets' Status column, and
# the rolled-up Overview sheet) and records the failure detail in the
# "Error Detail" column (L) for that row on each locale sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: the summary row for 26c9a277-...md now reflects the
#     failed status for both locale columns (B = zh-cn, C = de-de). ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: update Status (C3) and record the Error Detail (L3). ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("L3").Value = "Handback file name: jt5qlzlf.out is different with handoff file name: 26c9a277-901b-4ea2-b1f8-1505ff629220.4d9b6d90e775ed2ac02b10221c9985ac557198da.zh-cn."

# --- de-de sheet: update Status (C3) and record the Error Detail (L3). ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("L3").Value = "Handback file name: jt5qlzlf.out is different with handoff file name: 26c9a277-901b-4ea2-b1f8-1505ff629220.4d9b6d90e775ed2ac02b10221c9985ac557198da.de-de."
